$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1, matching the formatting of the
# existing header cells (e.g. G1 - bold, bordered, centered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add the corresponding data value in H2
$ws.Range("H2").Value = 0
